$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.112.56'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.14%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.819.78'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -0.42%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.07%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '312.16'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.21%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.000'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.02%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4468'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +5.36%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3747'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +2.02%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07496'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +3.73%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8730'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +3.08%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.95'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +1.19%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.806.30'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -1.12%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.736'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +1.02%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '94.31'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +5.27%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.348'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +1.02%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.07109'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +0.98%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.001'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -0.10%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008752'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -0.08%  '

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -0.16%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '15.01'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.82%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '27.127.77'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +0.00%  '

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +1.94%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.92'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +0.93%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.040.72'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -0.50%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.992'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +0.60%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.430'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +7.87%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '151.70'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +0.04%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.52'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +1.78%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.348'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +1.82%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '118.16'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +1.13%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08836'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +1.41%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.7664'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +3.88%  '

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -0.16%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.561'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +2.86%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.889'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -0.38%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.000'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -0.02%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.104'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +1.07%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01982'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +1.99%  '

$ws.Range("B39").Value = 'FraxShare'
$ws.Range("C39").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '7.481'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +2.15%  '

$ws.Range("B40").Value = 'Hedera'
$ws.Range("C40").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.05277'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +0.68%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.5318'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +4.76%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1720'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +1.90%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.851'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -0.82%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.197'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +11.27%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.714'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +1.68%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5053'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +6.61%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '10.64'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +1.33%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.702'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +3.07%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '105.64'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -0.08%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.9999'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -0.02%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06363'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +0.61%  '
